$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Target cluster = ECs (unchanged), update numeric values
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.574538
$ws.Range("H2").Value = 1.723614
$ws.Range("K2").Value = 2
$ws.Range("M2").Value = 33.6320075
$ws.Range("N2").Value = 67.264015
$ws.Range("O2").Value = 0.3908110491225105
$ws.Range("P2").Value = 0.3281482467916435
$ws.Range("Q2").Value = 19.322866325035
$ws.Range("R2").Value = 115.93719795021
$ws.Range("S2").Value = 0.3908110491225105
$ws.Range("T2").Value = 0.3281482467916435

# Row 3: Target cluster = FAPs (unchanged), update numeric values
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.574538
$ws.Range("H3").Value = 1.723614
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 12.015213
$ws.Range("N3").Value = 36.045639
$ws.Range("O3").Value = 0.1396193194224409
$ws.Range("P3").Value = 0.1758490515669409
$ws.Range("Q3").Value = 6.903196446594
$ws.Range("R3").Value = 62.128768019346
$ws.Range("S3").Value = 0.1396193194224409
$ws.Range("T3").Value = 0.1758490515669409

# Row 4: Target cluster becomes M1 (was Neutro), update numeric values
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.574538
$ws.Range("H4").Value = 1.723614
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.122567
$ws.Range("N4").Value = 0.367701
$ws.Range("O4").Value = 0.001424254495001488
$ws.Range("P4").Value = 0.001793833426290924
$ws.Range("Q4").Value = 0.070419399046
$ws.Range("R4").Value = 0.6337745914139999
$ws.Range("S4").Value = 0.001424254495001488
$ws.Range("T4").Value = 0.001793833426290924

# Row 5: Target cluster becomes M2 (was sCs), update numeric values
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.574538
$ws.Range("H5").Value = 1.723614
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.1323796666666667
$ws.Range("N5").Value = 0.397139
$ws.Range("O5").Value = 0.001538279759615546
$ws.Range("P5").Value = 0.001937447037358482
$ws.Range("Q5").Value = 0.07605714892733334
$ws.Range("R5").Value = 0.684514340346
$ws.Range("S5").Value = 0.001538279759615546
$ws.Range("T5").Value = 0.001937447037358482

# Row 6: new row - Target cluster = Neutro
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Wnt2"
$ws.Range("C6").Value = "Fzd4"
$ws.Range("D6").Value = "Neutro"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.574538
$ws.Range("H6").Value = 1.723614
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 1
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 20.59650833333334
$ws.Range("N6").Value = 61.789525
$ws.Range("O6").Value = 0.2393357883858265
$ws.Range("P6").Value = 0.3014408863169769
$ws.Range("Q6").Value = 11.83347670481667
$ws.Range("R6").Value = 106.50129034335
$ws.Range("S6").Value = 0.2393357883858265
$ws.Range("T6").Value = 0.3014408863169769

# Row 7: new row - Target cluster = sCs
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Wnt2"
$ws.Range("C7").Value = "Fzd4"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.574538
$ws.Range("H7").Value = 1.723614
$ws.Range("I7").Value = 1
$ws.Range("J7").Value = 1
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 19.558276
$ws.Range("N7").Value = 39.116552
$ws.Range("O7").Value = 0.2272713088146052
$ws.Range("P7").Value = 0.1908305348607893
$ws.Range("Q7").Value = 11.236972776488
$ws.Range("R7").Value = 67.421836658928
$ws.Range("S7").Value = 0.2272713088146052
$ws.Range("T7").Value = 0.1908305348607893
